$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44274
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 16000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 16000
$ws.Range("P2").Value = 889

# Row 3
$ws.Range("D3").Value = 44243
$ws.Range("J3").Value = 60

# Row 4
$ws.Range("D4").Value = 44260
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 16000
$ws.Range("P4").Value = 889

# Row 5
$ws.Range("D5").Value = 44252
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 18000
$ws.Range("P5").Value = 1000

# Row 6
$ws.Range("D6").Value = 44326
$ws.Range("I6").Value = "Especial"
$ws.Range("J6").Value = 15
$ws.Range("K6").Value = 18000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 18000
$ws.Range("P6").Value = 1000

# Row 7
$ws.Range("D7").Value = 44323
$ws.Range("I7").Value = "Especial"
$ws.Range("J7").Value = 70

# Row 8
$ws.Range("D8").Value = 44238
$ws.Range("J8").Value = 50

# Row 9
$ws.Range("D9").Value = 44245
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 18000
$ws.Range("P9").Value = 1000

# Row 10
$ws.Range("D10").Value = 44270
$ws.Range("I10").Value = "Primera"
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 16000
$ws.Range("P10").Value = 889

# Row 11
$ws.Range("D11").Value = 44365
$ws.Range("I11").Value = "Especial"
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 20000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 20000
$ws.Range("P11").Value = 1111

# Row 12
$ws.Range("D12").Value = 44232
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 22000
$ws.Range("L12").Value = 22000
$ws.Range("M12").Value = 22000
$ws.Range("P12").Value = 1222

# Row 13
$ws.Range("D13").Value = 44253
$ws.Range("J13").Value = 70

# Row 14
$ws.Range("D14").Value = 44246
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 60

# Row 15
$ws.Range("D15").Value = 44257
$ws.Range("I15").Value = "Primera"
$ws.Range("K15").Value = 16000
$ws.Range("L15").Value = 16000
$ws.Range("M15").Value = 16000
$ws.Range("P15").Value = 889

# Row 16
$ws.Range("D16").Value = 44285
$ws.Range("K16").Value = 18000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 18000
$ws.Range("P16").Value = 1000

# Row 17
$ws.Range("D17").Value = 44371
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 20000
$ws.Range("P17").Value = 1111

# Row 18
$ws.Range("D18").Value = 44302
$ws.Range("I18").Value = "Especial"
$ws.Range("J18").Value = 70

# Row 19
$ws.Range("D19").Value = 44242
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 20000
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = 20000
$ws.Range("P19").Value = 1111

# Row 20
$ws.Range("D20").Value = 44271

# Row 22
$ws.Range("D22").Value = 44259
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 70
$ws.Range("K22").Value = 16000
$ws.Range("L22").Value = 16000
$ws.Range("M22").Value = 16000
$ws.Range("P22").Value = 889

# Row 23
$ws.Range("D23").Value = 44309
$ws.Range("J23").Value = 80

# Row 24
$ws.Range("D24").Value = 44350
$ws.Range("J24").Value = 20
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 20000
$ws.Range("P24").Value = 1111

# Row 25
$ws.Range("D25").Value = 44267
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 70
$ws.Range("K25").Value = 16000
$ws.Range("L25").Value = 16000
$ws.Range("M25").Value = 16000
$ws.Range("P25").Value = 889

# Row 26
$ws.Range("D26").Value = 44292
$ws.Range("K26").Value = 17000
$ws.Range("L26").Value = 17000
$ws.Range("M26").Value = 17000
$ws.Range("P26").Value = 944

# Row 27
$ws.Range("D27").Value = 44299
$ws.Range("I27").Value = "Especial"
$ws.Range("J27").Value = 80

# Row 28
$ws.Range("D28").Value = 44357
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = 20000
$ws.Range("L28").Value = 20000
$ws.Range("M28").Value = 20000
$ws.Range("P28").Value = 1111

# Row 29
$ws.Range("D29").Value = 44320
$ws.Range("I29").Value = "Especial"
$ws.Range("J29").Value = 90
$ws.Range("K29").Value = 17500
$ws.Range("L29").Value = 17500
$ws.Range("M29").Value = 17500
$ws.Range("P29").Value = 972

# Row 30
$ws.Range("D30").Value = 44306
$ws.Range("I30").Value = "Especial"
$ws.Range("J30").Value = 80

# Row 31
$ws.Range("D31").Value = 44251
$ws.Range("J31").Value = 20
$ws.Range("K31").Value = 18000
$ws.Range("L31").Value = 18000
$ws.Range("M31").Value = 18000
$ws.Range("P31").Value = 1000

# Row 32
$ws.Range("D32").Value = 44295
$ws.Range("J32").Value = 80
$ws.Range("K32").Value = 16000
$ws.Range("L32").Value = 16000
$ws.Range("M32").Value = 16000
$ws.Range("P32").Value = 889

# Row 33
$ws.Range("D33").Value = 44316

# Row 34
$ws.Range("D34").Value = 44313
$ws.Range("I34").Value = "Especial"
$ws.Range("J34").Value = 80
$ws.Range("K34").Value = 18000
$ws.Range("L34").Value = 18000
$ws.Range("M34").Value = 18000
$ws.Range("P34").Value = 1000

# Row 35
$ws.Range("D35").Value = 44280
$ws.Range("J35").Value = 40
$ws.Range("K35").Value = 18000
$ws.Range("L35").Value = 18000
$ws.Range("M35").Value = 18000
$ws.Range("P35").Value = 1000

# Row 36
$ws.Range("D36").Value = 44239
$ws.Range("I36").Value = "Especial"
$ws.Range("K36").Value = 20000
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = 20000
$ws.Range("P36").Value = 1111

# Row 37
$ws.Range("D37").Value = 44250
$ws.Range("J37").Value = 60
$ws.Range("K37").Value = 18000
$ws.Range("L37").Value = 18000
$ws.Range("M37").Value = 18000
$ws.Range("P37").Value = 1000

# Row 38
$ws.Range("D38").Value = 44301
$ws.Range("J38").Value = 30

# Row 39
$ws.Range("D39").Value = 44278
$ws.Range("J39").Value = 70
$ws.Range("K39").Value = 18000
$ws.Range("L39").Value = 18000
$ws.Range("M39").Value = 18000
$ws.Range("P39").Value = 1000

# Row 40
$ws.Range("D40").Value = 44236
$ws.Range("J40").Value = 60

# Row 41
$ws.Range("D41").Value = 44264
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 80
$ws.Range("K41").Value = 16000
$ws.Range("L41").Value = 16000
$ws.Range("M41").Value = 16000
$ws.Range("P41").Value = 889

# Row 42
$ws.Range("D42").Value = 44327
$ws.Range("J42").Value = 80
$ws.Range("K42").Value = 16000
$ws.Range("L42").Value = 16000
$ws.Range("M42").Value = 16000
$ws.Range("P42").Value = 889

# Row 43
$ws.Range("D43").Value = 44312
$ws.Range("J43").Value = 15
